# Append the "04.03" row (row 12) to the daily COVID summary sheet,
# matching the other existing date rows (A2:A11) in both data and style.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A holds date labels ("23.02", "24.02", ... "03.03") stored as
# plain text. "04.03" looks like a number to Excel, so a plain
# `Value = "04.03"` assignment would be auto-coerced to the number 4.03.
# Prefixing with a single quote forces Excel to keep it as literal text,
# exactly like the other labels in the column.
$ws.Range("A12").Value = "'04.03"

# A12 should look just like the other date-label cells above it (bold,
# centered, bordered "style 1"). Copy that formatting down from A11.
$ws.Range("A11").Copy()
$ws.Range("A12").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$excel.CutCopyMode = $false

# Fill in the day's figures.
$ws.Range("B12").Value = 1346
$ws.Range("C12").Value = 295
$ws.Range("D12").Value = 1065
$ws.Range("E12").Value = 2706
$ws.Range("F12").Value = 276
$ws.Range("G12").Value = 107
$ws.Range("H12").Value = 3089
$ws.Range("I12").Value = 29837
